$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '42.918.42'
Set-TextValue 'E2' '  +1.03%  '
Set-TextValue 'D3' '2.290.47'
Set-TextValue 'E3' '  -0.49%  '
Set-TextValue 'E4' '  -0.42%  '
Set-TextValue 'D5' '314.20'
Set-TextValue 'E5' '  -0.63%  '
Set-TextValue 'D6' '105.31'
Set-TextValue 'E6' '  +1.72%  '
Set-TextValue 'E7' '  -0.47%  '
Set-TextValue 'E8' '  -0.09%  '
Set-TextValue 'D9' '0.606'
Set-TextValue 'E9' '  -0.25%  '
Set-TextValue 'D10' '39.74'
Set-TextValue 'E10' '  +0.08%  '
Set-TextValue 'D11' '0.0906'
Set-TextValue 'E11' '  -0.41%  '
Set-TextValue 'D12' '8.42'
Set-TextValue 'E12' '  +0.67%  '
Set-TextValue 'E13' '  +2.63%  '
Set-TextValue 'D14' '0.999'
Set-TextValue 'E14' '  +3.78%  '
Set-TextValue 'D15' '15.30'
Set-TextValue 'E15' '  +0.26%  '
Set-TextValue 'D16' '2.638.88'
Set-TextValue 'E16' '  -0.47%  '
Set-TextValue 'D17' '2.286.92'
Set-TextValue 'E17' '  -0.71%  '
Set-TextValue 'D18' '42.817.00'
Set-TextValue 'E18' '  +0.89%  '
Set-TextValue 'D19' '7.41'
Set-TextValue 'E19' '  -0.69%  '
Set-TextValue 'D20' '13.79'
Set-TextValue 'E20' '  +22.28%  '
Set-TextValue 'D21' '0.0000106'
Set-TextValue 'E21' '  -0.27%  '
Set-TextValue 'D22' '73.97'
Set-TextValue 'E22' '  +0.77%  '
Set-TextValue 'D23' '3.54'
Set-TextValue 'E23' '  +0.16%  '
Set-TextValue 'D24' '265.66'
Set-TextValue 'E24' '  -3.85%  '
Set-TextValue 'D25' '2.22'
Set-TextValue 'E25' '  -2.81%  '
Set-TextValue 'E26' '  +0.49%  '
Set-TextValue 'B27' 'Filecoin'
Set-TextValue 'C27' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D27' '7.35'
Set-TextValue 'E27' '  +25.17%  '
Set-TextValue 'B28' 'Cosmos'
Set-TextValue 'C28' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D28' '10.89'
Set-TextValue 'E28' '  +0.40%  '
Set-TextValue 'D29' '2.35'
Set-TextValue 'E29' '  -0.08%  '
Set-TextValue 'D30' '22.51'
Set-TextValue 'E30' '  -1.12%  '
Set-TextValue 'D31' '37.48'
Set-TextValue 'E31' '  +1.79%  '
Set-TextValue 'D32' '167.15'
Set-TextValue 'E32' '  +1.00%  '
Set-TextValue 'D33' '0.0877'
Set-TextValue 'E33' '  +0.34%  '
Set-TextValue 'E34' '  -2.62%  '
Set-TextValue 'E35' '  -1.14%  '
Set-TextValue 'E36' '  -3.81%  '
Set-TextValue 'D37' '4.56'
Set-TextValue 'E37' '  -0.16%  '
Set-TextValue 'E38' '  -4.07%  '
Set-TextValue 'D39' '3.82'
Set-TextValue 'E39' '  +3.07%  '
Set-TextValue 'E40' '  -3.66%  '
Set-TextValue 'E41' '  +4.81%  '
Set-TextValue 'D42' '70.97'
Set-TextValue 'E42' '  +1.88%  '
Set-TextValue 'D43' '0.232'
Set-TextValue 'E43' '  +2.43%  '
Set-TextValue 'D44' '94.38'
Set-TextValue 'E44' '  -0.35%  '
Set-TextValue 'E45' '  +0.06%  '
Set-TextValue 'D46' '12.23'
Set-TextValue 'E46' '  +1.37%  '
Set-TextValue 'D47' '1.736.07'
Set-TextValue 'E47' '  +9.17%  '
Set-TextValue 'D48' '113.72'
Set-TextValue 'E48' '  +0.60%  '
Set-TextValue 'D49' '79.71'
Set-TextValue 'E49' '  -1.80%  '
Set-TextValue 'D50' '8.75'
Set-TextValue 'E50' '  -2.38%  '
Set-TextValue 'D51' '5.18'
Set-TextValue 'E51' '  -0.57%  '

Write-Host "Applied cryptos update."
